$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$korean = @(
  "안녕하세요. 이메일 잘 받았습니다.",
  "회의 일정은 언제인가요?",
  "보고서는 내일까지 제출해주세요.",
  "결재를 부탁드립니다.",
  "이번 프로젝트의 목표는 무엇인가요?",
  "예산안을 검토해주세요.",
  "다음 주 출장 예정입니다.",
  "회의실 예약했습니다.",
  "잠시 통화 가능하신가요?",
  "지금은 회의 중입니다.",
  "나중에 다시 연락드리겠습니다.",
  "담당자 연결해 드리겠습니다.",
  "협조해주셔서 감사합니다.",
  "문의사항이 있으시면 연락주세요.",
  "첨부 파일을 확인해주세요.",
  "빠른 회신 부탁드립니다.",
  "거래처 미팅이 잡혔습니다.",
  "계약서 초안을 보냈습니다.",
  "수정 사항이 있습니다.",
  "승인되었습니다.",
  "반려되었습니다.",
  "추가 정보가 필요합니다.",
  "일정을 변경해야 할 것 같습니다.",
  "마감 기한을 지켜주세요.",
  "진행 상황을 공유해주세요.",
  "문제가 발생했습니다.",
  "해결 방안을 모색 중입니다.",
  "성공적으로 완료되었습니다.",
  "실적 보고회가 있습니다.",
  "야근을 해야 합니다.",
  "주말 근무는 없습니다.",
  "연차 휴가를 신청합니다.",
  "병가를 냈습니다.",
  "조퇴하겠습니다.",
  "지각해서 죄송합니다.",
  "출근했습니다.",
  "퇴근하겠습니다.",
  "점심 식사 맛있게 하세요.",
  "오늘 회식 참석하시나요?",
  "명함 좀 주시겠습니까?",
  "제 명함입니다.",
  "소개해 드리겠습니다.",
  "우리 팀장님입니다.",
  "신입 사원입니다.",
  "인수인계 중입니다.",
  "업무 분장이 필요합니다.",
  "효율성을 높여야 합니다.",
  "비용을 절감해야 합니다.",
  "고객 만족도가 중요합니다.",
  "시장 조사를 했습니다.",
  "경쟁사를 분석했습니다.",
  "전략을 세웠습니다.",
  "목표를 달성했습니다.",
  "성과급이 지급됩니다.",
  "승진을 축하합니다.",
  "퇴사하게 되었습니다.",
  "그동안 감사했습니다.",
  "송별회가 있습니다.",
  "환영회가 있습니다.",
  "워크숍을 갑니다.",
  "세미나에 참석합니다.",
  "교육을 받습니다.",
  "자격증을 취득했습니다.",
  "외국어 능력이 필요합니다.",
  "프레젠테이션을 잘했습니다.",
  "질문 있습니까?",
  "이해가 되시나요?",
  "다시 설명해 드리겠습니다.",
  "요점을 정리해주세요.",
  "회의록을 작성해주세요.",
  "아이디어가 좋습니다.",
  "다른 의견 있습니까?",
  "만장일치로 통과되었습니다.",
  "반대 의견이 있습니다.",
  "절충안을 찾읍시다.",
  "결론을 내립시다.",
  "다음 안건으로 넘어갑시다.",
  "시간이 부족합니다.",
  "잠시 쉬었다 합시다.",
  "오늘 회의는 여기까지입니다.",
  "수고 많으셨습니다.",
  "좋은 주말 보내세요.",
  "월요일에 뵙겠습니다.",
  "복사가 안 됩니다.",
  "프린터 용지가 없습니다.",
  "인터넷이 느립니다.",
  "서버가 다운되었습니다.",
  "보안이 중요합니다.",
  "비밀번호를 변경하세요.",
  "로그인이 안 됩니다.",
  "IT 팀에 문의하세요.",
  "탕비실에 커피가 없습니다.",
  "에어컨이 너무 춥습니다.",
  "히터 좀 틀어주세요.",
  "창문 좀 열어주세요.",
  "청소 업체 불라주세요.",
  "택배가 도착했습니다.",
  "우편물 확인해주세요.",
  "주차권 필요하세요?",
  "엘리베이터가 고장났습니다."
)

$english = @(
  "Hello. I received your email.",
  "When is the meeting schedule?",
  "Please submit the report by tomorrow.",
  "Please approve this.",
  "What is the goal of this project?",
  "Please review the budget proposal.",
  "I am scheduled for a business trip next week.",
  "I booked a meeting room.",
  "Are you available for a call?",
  "I am in a meeting right now.",
  "I will contact you later.",
  "I will connect you to the person in charge.",
  "Thank you for your cooperation.",
  "Please contact me if you have any questions.",
  "Please check the attached file.",
  "Please reply as soon as possible.",
  "A meeting with the client is scheduled.",
  "I sent the draft of the contract.",
  "There are some revisions.",
  "It has been approved.",
  "It has been rejected.",
  "More information is needed.",
  "I think we need to reschedule.",
  "Please meet the deadline.",
  "Please share the progress.",
  "A problem has occurred.",
  "We are looking for a solution.",
  "It has been successfully completed.",
  "There is a performance report meeting.",
  "I have to work overtime.",
  "There is no weekend work.",
  "I apply for annual leave.",
  "I took sick leave.",
  "I will leave early.",
  "I am sorry for being late.",
  "I have arrived at work.",
  "I am leaving work.",
  "Enjoy your lunch.",
  "Are you attending the company dinner today?",
  "Could you give me your business card?",
  "Here is my business card.",
  "Let me introduce you.",
  "This is our team leader.",
  "This is a new employee.",
  "I am in the middle of handover.",
  "Division of duties is needed.",
  "We need to increase efficiency.",
  "We need to cut costs.",
  "Customer satisfaction is important.",
  "We conducted market research.",
  "We analyzed competitors.",
  "We established a strategy.",
  "We achieved the goal.",
  "Performance bonuses will be paid.",
  "Congratulations on your promotion.",
  "I am resigning.",
  "Thank you for everything.",
  "There is a farewell party.",
  "There is a welcome party.",
  "We are going to a workshop.",
  "I am attending a seminar.",
  "I am receiving training.",
  "I obtained a certification.",
  "Foreign language skills are required.",
  "You did a good job on the presentation.",
  "Do you have any questions?",
  "Do you understand?",
  "I will explain it again.",
  "Please summarize the main points.",
  "Please write the meeting minutes.",
  "That's a good idea.",
  "Any other opinions?",
  "It was passed unanimously.",
  "There is a dissenting opinion.",
  "Let's find a compromise.",
  "Let's draw a conclusion.",
  "Let's move on to the next item.",
  "We are running out of time.",
  "Let's take a short break.",
  "That's all for today's meeting.",
  "Thank you for your hard work.",
  "Have a good weekend.",
  "See you on Monday.",
  "The copier is not working.",
  "The printer is out of paper.",
  "The internet is slow.",
  "The server is down.",
  "Security is important.",
  "Please change your password.",
  "I can't log in.",
  "Please contact the IT team.",
  "There is no coffee in the pantry.",
  "The air conditioner is too cold.",
  "Please turn on the heater.",
  "Please open the window.",
  "Please call the cleaning service.",
  "A package has arrived.",
  "Please check the mail.",
  "Do you need a parking ticket?",
  "The elevator is broken."
)

for ($i = 0; $i -lt $korean.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $korean[$i]
    $ws.Cells.Item($row, 2).Value = $english[$i]
}

